$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 145; A = "enior Backend Engineer (Go, Java/Scala, AWS, REST APIs)"; B = "https://www.dice.com/job-detail/826ae36d-1a8d-4076-8467-52c806376af3?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"; C = "Hybrid in McLean, Virginia"; D = "Contract"; E = "`$60 - `$70"; F = "InfiCare Technologies" },
    @{ Row = 146; A = "Go Developer"; B = "https://www.dice.com/job-detail/e3755cd6-3b09-4987-b2a6-0c95e57fab4c?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"; C = "Alpharetta, Georgia"; D = "Contract, Third Party"; E = "Depends on Experience"; F = "Astir IT Solutions" },
    @{ Row = 147; A = "Full Stack AWS Cloud Engineer - W2 - onsite, must be local to Chicago, IL (Posted by SAM)"; B = "https://www.dice.com/job-detail/759fca36-016d-4fc9-ac21-36a2eef94a10?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"; C = "Chicago, Illinois"; D = "Contract"; E = "Depends on Experience"; F = "NimbusAITech LLC" },
    @{ Row = 148; A = "SAP PP Consultant || Onsite in Normal/Chicago, IL (Relocation will work) || Automotive industry exp."; B = "https://www.dice.com/job-detail/1525e71a-a575-4d66-96a2-d6cc955fcda7?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"; C = "Texas"; D = "Contract, Third Party"; E = "Depends on Experience"; F = "Oraapps Inc" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
